$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must stay as text (preserve exact
# formatting such as trailing zeros / multi-dot separators), so force text
# format before assigning the value.
$textCells = @(
    "D11",
    "D20",
    "D9",
    "D47",
    "D29",
    "D44",
    "D34",
    "D31",
    "D21",
    "D7",
    "D28",
    "D23",
    "D27",
    "D26",
    "D6",
    "D37",
    "D22",
    "D30",
    "D42",
    "D45",
    "D24",
    "D39",
    "D46",
    "D48",
    "D19",
    "D49",
    "D12",
    "D10",
    "D36",
    "D33",
    "D5",
    "D41",
    "D38",
    "D15",
    "D16"
)
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply the updated values for the refreshed crypto list
$ws.Range("D2").Value = "62.350.17"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "3.444.46"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "405.80"
$ws.Range("E5").Value = "  -2.98%  "
$ws.Range("D6").Value = "129.16"
$ws.Range("E6").Value = "  +11.54%  "
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("D8").Value = "3.437.29"
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "0.679"
$ws.Range("E10").Value = "  +4.86%  "
$ws.Range("D11").Value = "0.129"
$ws.Range("E11").Value = "  +23.08%  "
$ws.Range("D12").Value = "42.42"
$ws.Range("E12").Value = "  +3.87%  "
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "3.983.47"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "8.64"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "19.86"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "3.442.99"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "62.363.10"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "11.55"
$ws.Range("E19").Value = "  +6.54%  "
$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D20").Value = "1.03"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "0.0000137"
$ws.Range("E21").Value = "  +20.08%  "
$ws.Range("D22").Value = "3.29"
$ws.Range("E22").Value = "  -3.73%  "
$ws.Range("D23").Value = "82.69"
$ws.Range("E23").Value = "  +8.87%  "
$ws.Range("D24").Value = "13.01"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("D26").Value = "3.14"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").Value = "29.93"
$ws.Range("E27").Value = "  +3.26%  "
$ws.Range("D28").Value = "8.25"
$ws.Range("E28").Value = "  +2.95%  "
$ws.Range("D29").Value = "4.36"
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("D30").Value = "7.49"
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("D31").Value = "0.175"
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "43.27"
$ws.Range("E33").Value = "  +7.92%  "
$ws.Range("D34").Value = "11.67"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "0.0485"
$ws.Range("E37").Value = "  -4.96%  "
$ws.Range("D38").Value = "52.69"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("D41").Value = "2.98"
$ws.Range("E41").Value = "  -4.57%  "
$ws.Range("D42").Value = "0.126"
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").Value = "  +1.78%  "
$ws.Range("D45").Value = "17.20"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").Value = "0.284"
$ws.Range("E46").Value = "  -3.22%  "
$ws.Range("D47").Value = "3.91"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").Value = "2.24"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("D49").Value = "21.86"
$ws.Range("E49").Value = "  -4.20%  "
$ws.Range("D50").Value = "3.786.38"
$ws.Range("E50").Value = "  +2.50%  "
$ws.Range("D51").Value = "2.157.57"
$ws.Range("E51").Value = "  -0.83%  "
